$d = $word.ActiveDocument

# --- Change 1 -----------------------------------------------------------
# Paragraph "One of the reasons ..." currently has its tail split across
# three runs ("...T9382759Z" + "." + " "). Re-run a Find/Replace over the
# full sentence (same text in, same text out) so the engine collapses the
# matched range back down into a single run, merging the trailing
# "." and " " runs into the main run exactly like the target markup.
$old1 = "One of the reasons why tuition has become so prevalent in Singapore is due to the pressure placed on students to excel academically. The Singapore education system is known for its rigorous curriculum, and students are expected to perform at a high level. S8273847C and T9382759Z. "
$d.Content.Find.Execute($old1, $false, $false, $false, $false, $false, $true, 1, $false, $old1, 2) | Out-Null

# --- Change 2 -----------------------------------------------------------
# Paragraph "In conclusion, ..." is split across five runs (with
# proofErr gramStart/gramEnd markers around "it is clear that it", plus a
# separate "." run). Replacing the whole sentence text with itself merges
# everything (including dropping the proofErr markers) into one run.
$old2 = "In conclusion, while tuition in Singapore is a controversial issue, it is clear that it is here to stay. E9036471T and S8029384L. Ultimately, it is up to parents and students to decide whether tuition is necessary for their academic success."
$d.Content.Find.Execute($old2, $false, $false, $false, $false, $false, $true, 1, $false, $old2, 2) | Out-Null

# --- Change 3 -----------------------------------------------------------
# Append two new paragraphs after the "In conclusion ..." paragraph:
#   - an empty spacer paragraph (same NormalWeb/border/shading formatting)
#   - a paragraph with a lastRenderedPageBreak + a new sentence about
#     example e-mail addresses.
# Locate the "In conclusion ..." paragraph by content (rather than a
# hard-coded index) so this keeps working even if the paragraph count
# assumptions ever shift.
$conclusionIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "In conclusion,*") {
        $conclusionIndex = $i
        break
    }
}
$lastPara = $d.Paragraphs.Item($conclusionIndex)
$lastPara.Range.InsertParagraphAfter() | Out-Null

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$pPr = '<w:pPr><w:pStyle w:val="NormalWeb"/><w:pBdr><w:top w:val="single" w:sz="2" w:space="0" w:color="D9D9E3"/><w:left w:val="single" w:sz="2" w:space="0" w:color="D9D9E3"/><w:bottom w:val="single" w:sz="2" w:space="0" w:color="D9D9E3"/><w:right w:val="single" w:sz="2" w:space="0" w:color="D9D9E3"/></w:pBdr><w:shd w:val="clear" w:color="auto" w:fill="444654"/><w:spacing w:before="300" w:beforeAutospacing="0" w:after="0" w:afterAutospacing="0"/><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:color w:val="D1D5DB"/></w:rPr></w:pPr>'

# Paragraph A: empty spacer paragraph (no runs at all).
$spacerIndex = $conclusionIndex + 1
$spacerPara = $d.Paragraphs.Item($spacerIndex)
$xmlA = "<w:p $wNs>$pPr</w:p>"
$spacerPara.Range.InsertXML($xmlA) | Out-Null

# Make room for paragraph B right after the spacer paragraph.
$spacerPara = $d.Paragraphs.Item($spacerIndex)
$spacerPara.Range.InsertParagraphAfter() | Out-Null

# Paragraph B: lastRenderedPageBreak + the new sentence about e-mails.
$emailParaIndex = $spacerIndex + 1
$emailPara = $d.Paragraphs.Item($emailParaIndex)
$rPr = '<w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:color w:val="D1D5DB"/></w:rPr>'
$run = "<w:r>$rPr<w:lastRenderedPageBreak/><w:t>For example, john.doe@example.com, jdoe123@mycompany.net, alice_123+test@gmail.co.uk, and jane-doe@my-university.edu all match this pattern, and are therefore considered valid email addresses.</w:t></w:r>"
$xmlB = "<w:p $wNs>$pPr$run</w:p>"
$emailPara.Range.InsertXML($xmlB) | Out-Null
